$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinbase API access secured: update G8, H8, I8, and J2 from "No" to "Yes"
$ws.Range("G8").Value = "Yes"
$ws.Range("H8").Value = "Yes"
$ws.Range("I8").Value = "Yes"
$ws.Range("J2").Value = "Yes"

# Update the selected cell to J8, matching the saved selection state
$ws.Range("J8").Select()
